$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 48.091872
$ws.Range("H2").Value2 = 144.275616
$ws.Range("I2").Value2 = 0.421093842675958
$ws.Range("J2").Value2 = 0.423782205092405
$ws.Range("M2").Value2 = 28.85518433333334
$ws.Range("N2").Value2 = 86.56555300000001
$ws.Range("O2").Value2 = 0.1999651185353207
$ws.Range("P2").Value2 = 0.2044513327926365
$ws.Range("Q2").Value2 = 1387.699831495072
$ws.Range("R2").Value2 = 12489.29848345565
$ws.Range("S2").Value2 = 0.08420408016519162
$ws.Range("T2").Value2 = 0.08664283664494463
$ws.Range("G3").Value2 = 48.091872
$ws.Range("H3").Value2 = 144.275616
$ws.Range("I3").Value2 = 0.421093842675958
$ws.Range("J3").Value2 = 0.423782205092405
$ws.Range("O3").Value2 = 0.3546352265743414
$ws.Range("P3").Value2 = 0.3625914622481308
$ws.Range("Q3").Value2 = 2461.065448634752
$ws.Range("R3").Value2 = 22149.58903771277
$ws.Range("S3").Value2 = 0.1493347103064485
$ws.Range("T3").Value2 = 0.1536598094191924
$ws.Range("G4").Value2 = 48.091872
$ws.Range("H4").Value2 = 144.275616
$ws.Range("I4").Value2 = 0.421093842675958
$ws.Range("J4").Value2 = 0.423782205092405
$ws.Range("M4").Value2 = 29.393479
$ws.Range("N4").Value2 = 88.180437
$ws.Range("O4").Value2 = 0.2036954761578358
$ws.Range("P4").Value2 = 0.2082653809291453
$ws.Range("Q4").Value2 = 1413.587429702688
$ws.Range("R4").Value2 = 12722.28686732419
$ws.Range("S4").Value2 = 0.08577491079101209
$ws.Range("T4").Value2 = 0.08825916237456288
$ws.Range("G5").Value2 = 48.091872
$ws.Range("H5").Value2 = 144.275616
$ws.Range("I5").Value2 = 0.421093842675958
$ws.Range("J5").Value2 = 0.423782205092405
$ws.Range("M5").Value2 = 9.499066500000001
$ws.Range("N5").Value2 = 18.998133
$ws.Range("O5").Value2 = 0.0658280999596015
$ws.Range("P5").Value2 = 0.04486996822421697
$ws.Range("Q5").Value2 = 456.8278902374881
$ws.Range("R5").Value2 = 2740.967341424929
$ws.Range("S5").Value2 = 0.02771980756804567
$ws.Range("T5").Value2 = 0.01901509407648481
$ws.Range("G6").Value2 = 48.091872
$ws.Range("H6").Value2 = 144.275616
$ws.Range("I6").Value2 = 0.421093842675958
$ws.Range("J6").Value2 = 0.423782205092405
$ws.Range("M6").Value2 = 25.37910966666666
$ws.Range("N6").Value2 = 76.13732899999999
$ws.Range("O6").Value2 = 0.1758760787729007
$ws.Range("P6").Value2 = 0.1798218558058706
$ws.Range("Q6").Value2 = 1220.528893563296
$ws.Range("R6").Value2 = 10984.76004206966
$ws.Range("S6").Value2 = 0.07406033384526023
$ws.Range("T6").Value2 = 0.07620530257722034
$ws.Range("I7").Value2 = 0.1230362686979479
$ws.Range("J7").Value2 = 0.1238217612582891
$ws.Range("M7").Value2 = 28.85518433333334
$ws.Range("N7").Value2 = 86.56555300000001
$ws.Range("O7").Value2 = 0.1999651185353207
$ws.Range("P7").Value2 = 0.2044513327926365
$ws.Range("Q7").Value2 = 405.4616620725832
$ws.Range("R7").Value2 = 3649.154958653248
$ws.Range("S7").Value2 = 0.02460296205432872
$ws.Range("T7").Value2 = 0.02531552411798884
$ws.Range("I8").Value2 = 0.1230362686979479
$ws.Range("J8").Value2 = 0.1238217612582891
$ws.Range("O8").Value2 = 0.3546352265743414
$ws.Range("P8").Value2 = 0.3625914622481308
$ws.Range("Q8").Value2 = 719.0803548615964
$ws.Range("R8").Value2 = 6471.723193754367
$ws.Range("S8").Value2 = 0.04363299502655832
$ws.Range("T8").Value2 = 0.04489671347278199
$ws.Range("I9").Value2 = 0.1230362686979479
$ws.Range("J9").Value2 = 0.1238217612582891
$ws.Range("M9").Value2 = 29.393479
$ws.Range("N9").Value2 = 88.180437
$ws.Range("O9").Value2 = 0.2036954761578358
$ws.Range("P9").Value2 = 0.2082653809291453
$ws.Range("Q9").Value2 = 413.0255662816213
$ws.Range("R9").Value2 = 3717.230096534592
$ws.Range("S9").Value2 = 0.02506193133711194
$ws.Range("T9").Value2 = 0.02578778627577525
$ws.Range("I10").Value2 = 0.1230362686979479
$ws.Range("J10").Value2 = 0.1238217612582891
$ws.Range("M10").Value2 = 9.499066500000001
$ws.Range("N10").Value2 = 18.998133
$ws.Range("O10").Value2 = 0.0658280999596015
$ws.Range("P10").Value2 = 0.04486996822421697
$ws.Range("Q10").Value2 = 133.477133493088
$ws.Range("R10").Value2 = 800.8628009585281
$ws.Range("S10").Value2 = 0.008099243794504905
$ws.Range("T10").Value2 = 0.00555587849312601
$ws.Range("I11").Value2 = 0.1230362686979479
$ws.Range("J11").Value2 = 0.1238217612582891
$ws.Range("M11").Value2 = 25.37910966666666
$ws.Range("N11").Value2 = 76.13732899999999
$ws.Range("O11").Value2 = 0.1758760787729007
$ws.Range("P11").Value2 = 0.1798218558058706
$ws.Range("Q11").Value2 = 356.6172327473848
$ws.Range("R11").Value2 = 3209.555094726464
$ws.Range("S11").Value2 = 0.02163913648544406
$ws.Range("T11").Value2 = 0.022265858898617
$ws.Range("G12").Value2 = 21.412221
$ws.Range("H12").Value2 = 64.23666299999999
$ws.Range("I12").Value2 = 0.1874860355013181
$ws.Range("J12").Value2 = 0.1886829905749125
$ws.Range("M12").Value2 = 28.85518433333334
$ws.Range("N12").Value2 = 86.56555300000001
$ws.Range("O12").Value2 = 0.1999651185353207
$ws.Range("P12").Value2 = 0.2044513327926365
$ws.Range("Q12").Value2 = 617.8535839410711
$ws.Range("R12").Value2 = 5560.682255469639
$ws.Range("S12").Value2 = 0.03749066731273841
$ws.Range("T12").Value2 = 0.03857648889834133
$ws.Range("G13").Value2 = 21.412221
$ws.Range("H13").Value2 = 64.23666299999999
$ws.Range("I13").Value2 = 0.1874860355013181
$ws.Range("J13").Value2 = 0.1886829905749125
$ws.Range("O13").Value2 = 0.3546352265743414
$ws.Range("P13").Value2 = 0.3625914622481308
$ws.Range("Q13").Value2 = 1095.754336234436
$ws.Range("R13").Value2 = 9861.789026109924
$ws.Range("S13").Value2 = 0.06648915267953495
$ws.Range("T13").Value2 = 0.0684148414539078
$ws.Range("G14").Value2 = 21.412221
$ws.Range("H14").Value2 = 64.23666299999999
$ws.Range("I14").Value2 = 0.1874860355013181
$ws.Range("J14").Value2 = 0.1886829905749125
$ws.Range("M14").Value2 = 29.393479
$ws.Range("N14").Value2 = 88.180437
$ws.Range("O14").Value2 = 0.2036954761578358
$ws.Range("P14").Value2 = 0.2082653809291453
$ws.Range("Q14").Value2 = 629.379668306859
$ws.Range("R14").Value2 = 5664.417014761731
$ws.Range("S14").Value2 = 0.0381900572743859
$ws.Range("T14").Value2 = 0.03929613490693448
$ws.Range("G15").Value2 = 21.412221
$ws.Range("H15").Value2 = 64.23666299999999
$ws.Range("I15").Value2 = 0.1874860355013181
$ws.Range("J15").Value2 = 0.1886829905749125
$ws.Range("M15").Value2 = 9.499066500000001
$ws.Range("N15").Value2 = 18.998133
$ws.Range("O15").Value2 = 0.0658280999596015
$ws.Range("P15").Value2 = 0.04486996822421697
$ws.Range("Q15").Value2 = 203.3961111916965
$ws.Range("R15").Value2 = 1220.376667150179
$ws.Range("S15").Value2 = 0.01234184948601016
$ws.Range("T15").Value2 = 0.008466199791546553
$ws.Range("G16").Value2 = 21.412221
$ws.Range("H16").Value2 = 64.23666299999999
$ws.Range("I16").Value2 = 0.1874860355013181
$ws.Range("J16").Value2 = 0.1886829905749125
$ws.Range("M16").Value2 = 25.37910966666666
$ws.Range("N16").Value2 = 76.13732899999999
$ws.Range("O16").Value2 = 0.1758760787729007
$ws.Range("P16").Value2 = 0.1798218558058706
$ws.Range("Q16").Value2 = 543.4231049659029
$ws.Range("R16").Value2 = 4890.807944693126
$ws.Range("S16").Value2 = 0.03297430874864866
$ws.Range("T16").Value2 = 0.03392932552418236
$ws.Range("G17").Value2 = 2.1734975
$ws.Range("H17").Value2 = 4.346995
$ws.Range("I17").Value2 = 0.01903120789977957
$ws.Range("J17").Value2 = 0.012768471746644
$ws.Range("M17").Value2 = 28.85518433333334
$ws.Range("N17").Value2 = 86.56555300000001
$ws.Range("O17").Value2 = 0.1999651185353207
$ws.Range("P17").Value2 = 0.2044513327926365
$ws.Range("Q17").Value2 = 62.71667101053917
$ws.Range("R17").Value2 = 376.300026063235
$ws.Range("S17").Value2 = 0.003805577743549753
$ws.Range("T17").Value2 = 0.002610531066326489
$ws.Range("G18").Value2 = 2.1734975
$ws.Range("H18").Value2 = 4.346995
$ws.Range("I18").Value2 = 0.01903120789977957
$ws.Range("J18").Value2 = 0.012768471746644
$ws.Range("O18").Value2 = 0.3546352265743414
$ws.Range("P18").Value2 = 0.3625914622481308
$ws.Range("Q18").Value2 = 111.2271029903767
$ws.Range("R18").Value2 = 667.36261794226
$ws.Range("S18").Value2 = 0.006749136725521725
$ws.Range("T18").Value2 = 0.004629738841289592
$ws.Range("G19").Value2 = 2.1734975
$ws.Range("H19").Value2 = 4.346995
$ws.Range("I19").Value2 = 0.01903120789977957
$ws.Range("J19").Value2 = 0.012768471746644
$ws.Range("M19").Value2 = 29.393479
$ws.Range("N19").Value2 = 88.180437
$ws.Range("O19").Value2 = 0.2036954761578358
$ws.Range("P19").Value2 = 0.2082653809291453
$ws.Range("Q19").Value2 = 63.88665312280249
$ws.Range("R19").Value2 = 383.3199187368149
$ws.Range("S19").Value2 = 0.003876570955004367
$ws.Range("T19").Value2 = 0.002659230632197841
$ws.Range("G20").Value2 = 2.1734975
$ws.Range("H20").Value2 = 4.346995
$ws.Range("I20").Value2 = 0.01903120789977957
$ws.Range("J20").Value2 = 0.012768471746644
$ws.Range("M20").Value2 = 9.499066500000001
$ws.Range("N20").Value2 = 18.998133
$ws.Range("O20").Value2 = 0.0658280999596015
$ws.Range("P20").Value2 = 0.04486996822421697
$ws.Range("Q20").Value2 = 20.64619729008375
$ws.Range("R20").Value2 = 82.58478916033501
$ws.Range("S20").Value2 = 0.001252788255978647
$ws.Range("T20").Value2 = 0.0005729209215437284
$ws.Range("G21").Value2 = 2.1734975
$ws.Range("H21").Value2 = 4.346995
$ws.Range("I21").Value2 = 0.01903120789977957
$ws.Range("J21").Value2 = 0.012768471746644
$ws.Range("M21").Value2 = 25.37910966666666
$ws.Range("N21").Value2 = 76.13732899999999
$ws.Range("O21").Value2 = 0.1758760787729007
$ws.Range("P21").Value2 = 0.1798218558058706
$ws.Range("Q21").Value2 = 55.16143141272583
$ws.Range("R21").Value2 = 330.968588476355
$ws.Range("S21").Value2 = 0.003347134219725081
$ws.Range("T21").Value2 = 0.00229605028528635
$ws.Range("G22").Value2 = 28.477822
$ws.Range("H22").Value2 = 85.433466
$ws.Range("I22").Value2 = 0.2493526452249964
$ws.Range("J22").Value2 = 0.2509445713277496
$ws.Range("M22").Value2 = 28.85518433333334
$ws.Range("N22").Value2 = 86.56555300000001
$ws.Range("O22").Value2 = 0.1999651185353207
$ws.Range("P22").Value2 = 0.2044513327926365
$ws.Range("Q22").Value2 = 821.7328032218554
$ws.Range("R22").Value2 = 7395.595228996698
$ws.Range("S22").Value2 = 0.04986183125951216
$ws.Range("T22").Value2 = 0.05130595206503523
$ws.Range("G23").Value2 = 28.477822
$ws.Range("H23").Value2 = 85.433466
$ws.Range("I23").Value2 = 0.2493526452249964
$ws.Range("J23").Value2 = 0.2509445713277496
$ws.Range("O23").Value2 = 0.3546352265743414
$ws.Range("P23").Value2 = 0.3625914622481308
$ws.Range("Q23").Value2 = 1457.331163498285
$ws.Range("R23").Value2 = 13115.98047148457
$ws.Range("S23").Value2 = 0.08842923183627796
$ws.Range("T23").Value2 = 0.09099035906095908
$ws.Range("G24").Value2 = 28.477822
$ws.Range("H24").Value2 = 85.433466
$ws.Range("I24").Value2 = 0.2493526452249964
$ws.Range("J24").Value2 = 0.2509445713277496
$ws.Range("M24").Value2 = 29.393479
$ws.Range("N24").Value2 = 88.180437
$ws.Range("O24").Value2 = 0.2036954761578358
$ws.Range("P24").Value2 = 0.2082653809291453
$ws.Range("Q24").Value2 = 837.062262922738
$ws.Range("R24").Value2 = 7533.560366304641
$ws.Range("S24").Value2 = 0.05079200580032155
$ws.Range("T24").Value2 = 0.05226306673967483
$ws.Range("G25").Value2 = 28.477822
$ws.Range("H25").Value2 = 85.433466
$ws.Range("I25").Value2 = 0.2493526452249964
$ws.Range("J25").Value2 = 0.2509445713277496
$ws.Range("M25").Value2 = 9.499066500000001
$ws.Range("N25").Value2 = 18.998133
$ws.Range("O25").Value2 = 0.0658280999596015
$ws.Range("P25").Value2 = 0.04486996822421697
$ws.Range("Q25").Value2 = 270.5127249531631
$ws.Range("R25").Value2 = 1623.076349718978
$ws.Range("S25").Value2 = 0.01641441085506211
$ws.Range("T25").Value2 = 0.01125987494151587
$ws.Range("G26").Value2 = 28.477822
$ws.Range("H26").Value2 = 85.433466
$ws.Range("I26").Value2 = 0.2493526452249964
$ws.Range("J26").Value2 = 0.2509445713277496
$ws.Range("M26").Value2 = 25.37910966666666
$ws.Range("N26").Value2 = 76.13732899999999
$ws.Range("O26").Value2 = 0.1758760787729007
$ws.Range("P26").Value2 = 0.1798218558058706
$ws.Range("Q26").Value2 = 722.7417676058126
$ws.Range("R26").Value2 = 6504.675908452313
$ws.Range("S26").Value2 = 0.04385516547382261
$ws.Range("T26").Value2 = 0.0451253185205646
